$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for time_taken column, copying the formatting from the
# existing "panel" header cell (E1) so it matches the other headers.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add time_taken values (as text) for each data row
$ws.Range("F2").Value = "2021-10-05 13:40:56.167939"
$ws.Range("F3").Value = "2021-10-05 13:40:56.167950"
$ws.Range("F4").Value = "2021-10-05 13:40:56.167954"
